$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the "duplicate_image_filename" column (E) with "NA" for rows 2-21,
# matching the other data columns already populated in that range.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
